$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 11: new data row
$ws.Range("A11").Value = 20230817
$ws.Range("B11").Value = 1
$ws.Range("C11").Value = 2
$ws.Range("D11").Value = 3
$ws.Range("E11").Value = 4
$ws.Range("F11").Value = 5
$ws.Range("G11").Value = 6

# Row 13: new data row
$ws.Range("A13").Value = 20230824
$ws.Range("B13").Value = 1
$ws.Range("C13").Value = 2
$ws.Range("D13").Value = 3
$ws.Range("E13").Value = 4
$ws.Range("F13").Value = 5
$ws.Range("G13").Value = 6

# Row 14: annotation row
$ws.Range("B14").Value = "7,14"
$ws.Range("C14").Value = 18
$ws.Range("D14").Value = "7,10"
$ws.Range("E14").Value = "good run"
$ws.Range("F14").Value = 19
$ws.Range("G14").Value = "really good run"

# Selection / view update
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 4
$ws.Range("H14").Select()
